$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B5 currently holds the text string "3". Copy it down to B6 first so the
# new row 6 keeps that same text-typed "3" (matching the moved annotation
# row), then turn B5 itself into a real number.
$ws.Range("B5").Copy($ws.Range("B6"))

# Row 5: change B5 from text "3" to numeric 3 (rest of row stays the same)
$ws.Range("B5").Value = 3

# Row 6: new row of data
$ws.Range("A6").Value = "Sunsi Wu"
$ws.Range("C6").Value = "无"
$ws.Range("D6").Value = "APC"
$ws.Range("E6").Value = "RES"
$ws.Range("F6").Value = "42be9703-0e9b-4ce8-962d-60bf1f233ce8"
$ws.Range("G6").Value = "SJCPLLpaW_annotated.xlsx"
$ws.Range("H6").Value = "The results show that DeePa achieves speedups compared to PyTorch and TensorFlow with all of the tested minibatch sizes."
